$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ======================================================================
# 1. Insert three new rows above the old row 10 ("Scaling" section).
#    This pushes the old rows 10-19 down to 13-22, matching the diff's
#    new row numbering.
# ======================================================================
$ws.Range("A10:A12").EntireRow.Insert()
$ws.Rows("10:12").RowHeight = 18.75

# ======================================================================
# 2. Row 9 gains two styled-but-empty cells (B9, C9), matching the
#    style used for the "value label" column elsewhere (copy from B8,
#    which carries the plain style=3 used on both B9 and C9).
# ======================================================================
$ws.Range("B8").Copy() | Out-Null
$ws.Range("B9").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("B9:C9").Value2 = ""

# ======================================================================
# 3. New row 10: "Mip Gap" sub-header, mirrors the "Power Scaling
#    Factor" sub-header (now at row 14) in style.
# ======================================================================
$ws.Range("B14:C14").Copy() | Out-Null
$ws.Range("B10:C10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("B10").Value2 = "Mip Gap"
$ws.Range("C10").Value2 = "[%]"

# ======================================================================
# 4. New row 11: the "pMIPGap" value row, mirrors the "pPowerScaling
#    Factor" value row (now at row 15) in style for B/C/E/F/H; the G
#    (unit) cell instead mirrors G8 ("Factor"-style single column, no
#    fill), matching the source row used in the diff.
# ======================================================================
$ws.Range("B15:H15").Copy() | Out-Null
$ws.Range("B11:H11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("G8").Copy() | Out-Null
$ws.Range("G11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("B11").Value2 = "pMIPGap"
$ws.Range("C11").Value2 = 0.05
$ws.Range("E11").Value2 = "Relative MIP gap"
$ws.Range("F11").Value2 = "The MIP solver will terminate (with an optimal result) when the gap between the lower and upper objective bound is less than pMIPGap"
$ws.Range("G11").Value2 = "Factor"
$ws.Range("H11").Value2 = 0.05

# ======================================================================
# 5. New row 12: blank spacer row, mirrors the blank spacer pattern
#    (e.g. now at row 16 / row 9).
# ======================================================================
$ws.Range("E16:H16").Copy() | Out-Null
$ws.Range("E12:H12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Re-set C5/C8 to force shared-string slot churn, matching the apparent
# target ordering where "No"/"gurobi" end up AFTER the new Mip-Gap strings.
$c5val = $ws.Range("C5").Value2
$c8val = $ws.Range("C8").Value2
$ws.Range("C5").Value2 = ""
$ws.Range("C8").Value2 = ""
$ws.Range("C5").Value2 = $c5val
$ws.Range("C8").Value2 = $c8val

Write-Host ("B1: " + $ws.Range("B1").Value2)
Write-Host ("B9: [" + $ws.Range("B9").Value2 + "]")
Write-Host ("B10: " + $ws.Range("B10").Value2)
Write-Host ("C10: " + $ws.Range("C10").Value2)
Write-Host ("B11: " + $ws.Range("B11").Value2)
Write-Host ("C11: " + $ws.Range("C11").Value2)
Write-Host ("E11: " + $ws.Range("E11").Value2)
Write-Host ("F11: " + $ws.Range("F11").Value2)
Write-Host ("G11: " + $ws.Range("G11").Value2)
Write-Host ("H11: " + $ws.Range("H11").Value2)
Write-Host ("B13 (Scaling header): " + $ws.Range("B13").Value2)
Write-Host ("B20 (General header): " + $ws.Range("B20").Value2)
Write-Host ("C22: " + $ws.Range("C22").Value2)
